# The deck's custom "Integral" design is switched to the built-in
# "Office Theme" design (the classic 12-slot Office colour scheme:
# black/white, dark/light neutrals and 6 accents, plus hyperlink colours).
# Re-create that swap through the theme-colour object model exposed on a
# slide (equivalent to picking "Office Theme" from the Design gallery,
# which rewrites the presentation's theme colours in place).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# ppThemeColorDark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $officeThemeRGB.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
